$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1811594202898551
$ws.Range("C2").Value = 0.5833333333333334
$ws.Range("J2").Value = 0.0108695652173913
$ws.Range("P2").Value = 0.1231884057971015
$ws.Range("S2").Value = 0.1014492753623188
$ws.Range("C3").Value = 0.02958579881656805
$ws.Range("J3").Value = 0.05917159763313609
$ws.Range("O3").Value = 0.005917159763313609
$ws.Range("P3").Value = 0.7100591715976331
$ws.Range("S3").Value = 0.1952662721893491
$ws.Range("J4").Value = 0.05263157894736842
$ws.Range("P4").Value = 0.5964912280701754
$ws.Range("S4").Value = 0.3508771929824561
$ws.Range("B6").Value = 0.06341463414634146
$ws.Range("D6").Value = 0.01463414634146342
$ws.Range("E6").Value = 0.00975609756097561
$ws.Range("F6").Value = 0.07804878048780488
$ws.Range("J6").Value = 0.375609756097561
$ws.Range("O6").Value = 0.01463414634146342
$ws.Range("Q6").Value = 0.1317073170731707
$ws.Range("R6").Value = 0.05365853658536585
$ws.Range("S6").Value = 0.2585365853658537
$ws.Range("B7").Value = 0.119496855345912
$ws.Range("D7").Value = 0.06289308176100629
$ws.Range("F7").Value = 0.05031446540880503
$ws.Range("J7").Value = 0.1509433962264151
$ws.Range("O7").Value = 0.01257861635220126
$ws.Range("Q7").Value = 0.1761006289308176
$ws.Range("R7").Value = 0.09433962264150944
$ws.Range("S7").Value = 0.3333333333333333
$ws.Range("B8").Value = 0.09090909090909091
$ws.Range("D8").Value = 0.01995565410199556
$ws.Range("F8").Value = 0.04212860310421286
$ws.Range("J8").Value = 0.1574279379157428
$ws.Range("O8").Value = 0.01995565410199556
$ws.Range("Q8").Value = 0.1507760532150776
$ws.Range("R8").Value = 0.09312638580931264
$ws.Range("S8").Value = 0.4257206208425721
$ws.Range("B9").Value = 0.0855614973262032
$ws.Range("D9").Value = 0.03208556149732621
$ws.Range("E9").Value = 0.0053475935828877
$ws.Range("F9").Value = 0.0267379679144385
$ws.Range("J9").Value = 0.1443850267379679
$ws.Range("O9").Value = 0.0053475935828877
$ws.Range("Q9").Value = 0.2299465240641711
$ws.Range("R9").Value = 0.09090909090909091
$ws.Range("S9").Value = 0.3796791443850268
$ws.Range("B10").Value = 0.1172893136403128
$ws.Range("D10").Value = 0.02780191138140747
$ws.Range("E10").Value = 0.0008688097306689834
$ws.Range("F10").Value = 0.06168549087749783
$ws.Range("J10").Value = 0.1277150304083406
$ws.Range("O10").Value = 0.01737619461337967
$ws.Range("Q10").Value = 0.2137271937445699
$ws.Range("R10").Value = 0.08079930495221546
$ws.Range("S10").Value = 0.3527367506516073
$ws.Range("G11").Value = 0.1428571428571428
$ws.Range("J11").Value = 0.1071428571428571
$ws.Range("K11").Value = 0.1830357142857143
$ws.Range("L11").Value = 0.5669642857142857
$ws.Range("G12").Value = 0.7938931297709924
$ws.Range("J12").Value = 0.1755725190839695
$ws.Range("L12").Value = 0.01526717557251908
$ws.Range("S12").Value = 0.01526717557251908
$ws.Range("G13").Value = 0.675
$ws.Range("J13").Value = 0.3
$ws.Range("S13").Value = 0.025
$ws.Range("F15").Value = 0.0154639175257732
$ws.Range("H15").Value = 0.2010309278350516
$ws.Range("I15").Value = 0.07731958762886598
$ws.Range("J15").Value = 0.2989690721649484
$ws.Range("K15").Value = 0.04639175257731959
$ws.Range("M15").Value = 0.0154639175257732
$ws.Range("O15").Value = 0.06701030927835051
$ws.Range("S15").Value = 0.2783505154639175
$ws.Range("F16").Value = 0.0160427807486631
$ws.Range("H16").Value = 0.1657754010695187
$ws.Range("I16").Value = 0.06417112299465241
$ws.Range("J16").Value = 0.5133689839572193
$ws.Range("K16").Value = 0.09090909090909091
$ws.Range("M16").Value = 0.0213903743315508
$ws.Range("O16").Value = 0.053475935828877
$ws.Range("S16").Value = 0.0748663101604278
$ws.Range("F17").Value = 0.05596107055961071
$ws.Range("H17").Value = 0.1873479318734793
$ws.Range("I17").Value = 0.07785888077858881
$ws.Range("J17").Value = 0.3673965936739659
$ws.Range("K17").Value = 0.06326034063260341
$ws.Range("M17").Value = 0.0194647201946472
$ws.Range("O17").Value = 0.07785888077858881
$ws.Range("S17").Value = 0.1508515815085158
$ws.Range("F18").Value = 0.03954802259887006
$ws.Range("H18").Value = 0.2203389830508475
$ws.Range("I18").Value = 0.0847457627118644
$ws.Range("J18").Value = 0.3559322033898305
$ws.Range("K18").Value = 0.0847457627118644
$ws.Range("M18").Value = 0.005649717514124294
$ws.Range("O18").Value = 0.05649717514124294
$ws.Range("S18").Value = 0.1525423728813559
$ws.Range("F19").Value = 0.02321428571428572
$ws.Range("H19").Value = 0.2366071428571428
$ws.Range("I19").Value = 0.1008928571428571
$ws.Range("J19").Value = 0.3392857142857143
$ws.Range("K19").Value = 0.1
$ws.Range("M19").Value = 0.02321428571428572
$ws.Range("N19").Value = 0.0008928571428571428
$ws.Range("O19").Value = 0.0625
$ws.Range("S19").Value = 0.1133928571428571
